# Generate Report for Handback
# Marks the two handed-off files as handed back (in sync with en-US),
# stamping the handback datetime and recording the target/handback
# file links for both the zh-cn and de-de localization sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Update the overall status text everywhere it's used (Overview
#    sheet + per-language Status columns all share this string).
# ---------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------
# 2. zh-cn sheet: stamp Latest Target File / Latest Handback File /
#    Latest Handback DateTime for both rows, and re-create the
#    hyperlinks collection so the new Target File links slot in
#    between the existing Source File Name hyperlinks.
# ---------------------------------------------------------------
$zhcn.Range("J2").Value = "3ea7805e-1905-438b-96d3-d506af6b7ad2.c07ecb2d0ee3c0284eef0d73ab4fa78b04d303b0.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-10-27 10:27:17"
$zhcn.Range("J3").Value = "fb9a253d-76aa-472d-ac07-f28f5a25a89c.80934fbed164e6f9e348594581e9fc4056c3f2b3.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-10-27 10:27:17"

# Only the new "Latest Target File" cells get hyperlinks; the existing
# Source File Name links in column A are left untouched.
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/215f6663d70064676b79cd28a720aeb8ab7c61f4/e2e/3ea7805e-1905-438b-96d3-d506af6b7ad2.md", "", "", "3ea7805e-1905-438b-96d3-d506af6b7ad2.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/215f6663d70064676b79cd28a720aeb8ab7c61f4/e2e/fb9a253d-76aa-472d-ac07-f28f5a25a89c.md", "", "", "fb9a253d-76aa-472d-ac07-f28f5a25a89c.md") | Out-Null

# ---------------------------------------------------------------
# 3. de-de sheet: same treatment, different timestamp + filenames.
# ---------------------------------------------------------------
$dede.Range("J2").Value = "3ea7805e-1905-438b-96d3-d506af6b7ad2.c07ecb2d0ee3c0284eef0d73ab4fa78b04d303b0.de-de.xlf"
$dede.Range("K2").Value = "2016-10-27 10:27:33"
$dede.Range("J3").Value = "fb9a253d-76aa-472d-ac07-f28f5a25a89c.80934fbed164e6f9e348594581e9fc4056c3f2b3.de-de.xlf"
$dede.Range("K3").Value = "2016-10-27 10:27:33"

$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/215f6663d70064676b79cd28a720aeb8ab7c61f4/e2e/3ea7805e-1905-438b-96d3-d506af6b7ad2.md", "", "", "3ea7805e-1905-438b-96d3-d506af6b7ad2.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/215f6663d70064676b79cd28a720aeb8ab7c61f4/e2e/fb9a253d-76aa-472d-ac07-f28f5a25a89c.md", "", "", "fb9a253d-76aa-472d-ac07-f28f5a25a89c.md") | Out-Null

# ---------------------------------------------------------------
# 4. Re-fit the columns that now hold longer text so the widths
#    track what Excel would compute after the content changes.
# ---------------------------------------------------------------
$overview.Columns("E").AutoFit() | Out-Null
$overview.Columns("F").AutoFit() | Out-Null
$zhcn.Columns("C").AutoFit() | Out-Null
$zhcn.Columns("I").AutoFit() | Out-Null
$zhcn.Columns("J").AutoFit() | Out-Null
$dede.Columns("C").AutoFit() | Out-Null
$dede.Columns("I").AutoFit() | Out-Null
$dede.Columns("J").AutoFit() | Out-Null
